$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$rng = $ws.Range("A2:N41")
$key = $ws.Range("D2:D41")
$rng.Sort($key, 1)
